# Actualización automática 2025-10-16 12:30:09
#
# Inserts a new sales-rep row ("ROCAFUERTE LOPEZ EVELYN ESTEFANIA") as row 16
# on both worksheets, pushing the existing row 16 ("VIEJO RIVAS MAYRA
# ANABELLE") and the trailing summary row down by one. The summary row's
# "0 de 15" / "4 de 15" / "2 de 15" labels on the first sheet become
# "0 de 16" / "4 de 16" / "2 de 16" to reflect the extra data row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new blank row above the current row 16; existing rows 16-17
# (and their formatting) shift down to 17-18.
$ws1.Rows.Item(16).Insert()

$ws1.Range("A16").Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws1.Range("B16").Value = "ROCAFUERTE LOPEZ EVELYN ESTEFANIA"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(16, $col).Value = 0
}

# Update the trailing "0 de 15" -> "0 de 16" summary row (now row 18).
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(18, $col)
    $cell.Value = ($cell.Value2 -replace "de 15", "de 16")
}

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(16).Insert()

$ws2.Range("A16").Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws2.Range("B16").Value = "ROCAFUERTE LOPEZ EVELYN ESTEFANIA"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(16, $col).Value = 0
}
